$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-11 03:26:47"
$wsZhCn.Range("G3").Value = "2016-01-11 03:27:37"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-11 03:27:00"
$wsDeDe.Range("G3").Value = "2016-01-11 03:27:58"
